$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (A:E) to match new report layout.
$ws.Columns.Item(1).ColumnWidth = 13.1666666667
$ws.Columns.Item(2).ColumnWidth = 20.1666666667
$ws.Columns.Item(3).ColumnWidth = 24.1666666667
$ws.Columns.Item(4).ColumnWidth = 11.1666666667
$ws.Columns.Item(5).ColumnWidth = 14.1666666667

# Full replacement data set (rows 2-25): ID_Categoria, Nome_Categoria, Produto, Quantidade, Data_Validade
$data = @(
    @(26, "LIMPEZA", "Sabaoembarra", 2, ""),
    @(24, "INFORMATICA", "Cabousb", 1, ""),
    @(20, "PAPELARIA", "Lapis", 5, ""),
    @(22, "MATERIAISSALADEAULA", "Apagador", 1, ""),
    @(24, "INFORMATICA", "Teste", 5, ""),
    @(27, "FERRAMENTAS", "Chavedefenda", 1, ""),
    @(20, "PAPELARIA", "Cadernoespiral", 50, ""),
    @(20, "PAPELARIA", "Borrachabranca", 20, ""),
    @(26, "LIMPEZA", "Detergente", 1, ""),
    @(26, "LIMPEZA", "Sabaoliquido", 200, "31/12/2025"),
    @(26, "LIMPEZA", "Cadernoespiral", 100, ""),
    @(26, "LIMPEZA", "Papela1Sulfite", 10, ""),
    @(26, "LIMPEZA", "Borrachabranca", 40, ""),
    @(20, "PAPELARIA", "Canetaesferograficaazul", 300, "31/12/2025"),
    @(24, "INFORMATICA", "Mouseusb", 1, ""),
    @(28, "MATERIALESCOLAR", "Cadernoespiral", 100, ""),
    @(24, "INFORMATICA", "Cabovga", 1, ""),
    @(20, "PAPELARIA", "Papela4Sulfite", 15, ""),
    @(28, "MATERIALESCOLAR", "Borrachabranca", 40, ""),
    @(20, "PAPELARIA", "Canetapreta", 4, "10/12/2025"),
    @(20, "PAPELARIA", "Canetaazul", 4, ""),
    @(26, "LIMPEZA", "Sabaoempo", 2, ""),
    @(24, "INFORMATICA", "Tecladousb", 0, ""),
    @(20, "PAPELARIA", "Canetavermelha", 100, "")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = 2 + $i
    $rowVals = $data[$i]
    $ws.Cells.Item($rowNum, 1).Value = $rowVals[0]
    $ws.Cells.Item($rowNum, 2).Value = $rowVals[1]
    $ws.Cells.Item($rowNum, 3).Value = $rowVals[2]
    $ws.Cells.Item($rowNum, 4).Value = $rowVals[3]
    if ($rowVals[4] -ne "") {
        # Force the date-like text to be stored as literal text (not
        # auto-converted to a date serial), then restore the default
        # "Normal" cell style so no stray number format sticks around.
        $ws.Cells.Item($rowNum, 5).NumberFormat = "@"
        $ws.Cells.Item($rowNum, 5).Value = $rowVals[4]
        $ws.Cells.Item($rowNum, 5).NumberFormat = "General"
        $ws.Cells.Item($rowNum, 5).Style = "Normal"
    }
}
